# Add materials for unit 7
#
# Applies the "Unit 7" (GradeRecorder / complex Firebase data) lesson block
# to the "2015" sheet, mirroring the target diff:
#   - New rows 78-85 added for the new "Unit 7" lesson block (names/times/notes)
#   - Row 67's C note is replaced ("(do Bolt...)" -> "(Skipped bolt)")
#   - Rows 65/66/68/69 get a day-of-week marker "M" in column C
#   - Rows 70/71 get day markers "T" in C; their old C notes move to D,
#     reworded for row 71
#   - Rows 72/73 get day markers "R" in C; row 73's old numeric note (11) is
#     replaced by "R" and its D note is reworded
#   - Row 74 gets a day marker "F" in column C
#
# New text values are written in the same order the original author entered
# them (new lesson rows/times first, then "6 videos", then the two rewritten
# class notes, then the "(Skipped bolt)" note, then the closing lab/total
# notes) so the rebuilt shared-string table lines up with the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- new rows for the Unit 7 lesson block ----
$ws.Range("A78").Value = "FirebaseComplexDataOverview"
$ws.Range("B78").Value = 0.20138888888888887
$ws.Range("B78").NumberFormat = "h:mm"

$ws.Range("A79").Value = "FirebaseDataDesign"
$ws.Range("B79").Value = 0.4694444444444445
$ws.Range("B79").NumberFormat = "h:mm"

$ws.Range("A80").Value = "GradeRecorderSetup"
$ws.Range("B80").Value = 0.53680555555555554
$ws.Range("B80").NumberFormat = "h:mm"

$ws.Range("A81").Value = "GradeRecorderFirebaseCode"
$ws.Range("B81").Value = 0.65416666666666667
$ws.Range("B81").NumberFormat = "h:mm"
$ws.Range("D81").Value = 44

$ws.Range("A82").Value = "GradeRecorderFirebaseCodePart2"
$ws.Range("B82").Value = 0.71319444444444446
$ws.Range("B82").NumberFormat = "h:mm"

$ws.Range("A83").Value = "LabGradeRecorder"
$ws.Range("B83").Value = 0.15833333333333333
$ws.Range("B83").NumberFormat = "h:mm"

$ws.Range("B85").Value = "6 videos"

# ---- rewritten class notes (new shared strings, in authoring order) ----
$ws.Range("D71").Value = "I dilly-dallied and it went ~40 min. Glad they did setup at home."
$ws.Range("D73").Value = "Short class. Plenty of time for Tyler to explain the system."
$ws.Range("C67").Value = "(Skipped bolt)"
$ws.Range("D83").Value = "21 then lab time."
$ws.Range("B84").Value = "65 min"

# ---- day-of-week markers in column C (reuse existing M/T/R/F shared strings) ----
$ws.Range("C65").Value = "M"
$ws.Range("C66").Value = "M"
$ws.Range("C68").Value = "M"
$ws.Range("C69").Value = "M"
$ws.Range("C70").Value = "T"
$ws.Range("D70").Value = "Do at home since SHA1 hash is painful."
$ws.Range("C71").Value = "T"
$ws.Range("C72").Value = "R"
$ws.Range("C73").Value = "R"
$ws.Range("C74").Value = "F"
$ws.Range("C78").Value = "T"
$ws.Range("C79").Value = "T"
$ws.Range("C80").Value = "T"
$ws.Range("C81").Value = "T"
$ws.Range("C82").Value = "R"
$ws.Range("C83").Value = "R"

# ---- the numeric note that used to live in C69 now lives in D69 ----
$ws.Range("D69").Value = 41

# ---- bold "Total" / summary row for the new block ----
$ws.Range("A84").Value = "Total"
$ws.Range("A84").Font.Bold = $true
$ws.Range("B84").Font.Bold = $true

# Move the visible selection/cursor to match the edited area.
[void]$ws.Range("B84").Select()
